$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B ("24<->31" / "All") changes from "Reserved" to a "Validity" bit-field
# like the other columns (C, D, E, F) - with 0=invalid / 1=valid meanings.
$ws.Range("B3").Value = "Validity"
$ws.Range("B4").Value = "0=invalid"
$ws.Range("B5").Value = "1=valid"

# The example bit value that used to live in B4 moves over to G4 (under the
# "8<->15" / "All" / "Reserved" column).
$ws.Range("G4").Value = 0

# Update the active selection to match the author's cursor position.
$ws.Range("H4").Select()
